$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Den Haag Bankaplein"
$ws.Cells.Item($row, 3).Value = "KDV"

# D69 holds a date-looking string that must stay plain text (as in the
# rest of the sheet), so force text format before assigning it, then
# drop back to the default "Normal" style so no stray style index is
# left attached to the cell (matches cells elsewhere in the sheet that
# carry no explicit style).
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-09-23"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
